$d = $word.ActiveDocument

# Full text of the target bullet, as it reads after the edit.
$full = "(BUG_FIX) Fix colliders on cat getting stuck at side of platforms and ground"

# Step 1: the original paragraph is made of two separate runs
# ("(BUG_FIX) " and "Fix colliders..."). Doing a self Find/Replace over the
# whole sentence collapses it back down into a single run, matching what
# Word does when you retype/reformat a whole selection.
$d.Content.Find.Execute($full, $true, $false, $false, $false, $false, `
                         $true, 1, $false, $full, 2) | Out-Null

# Step 2: find that paragraph and highlight it cyan, then split the run and
# drop the (document-unique) _GoBack bookmark into its new position. Adding
# a bookmark with a name that already exists elsewhere in the document moves
# it here, which is exactly what removes it from the old "Quit Button" spot.
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "*Fix colliders on cat getting stuck*") {

        # Highlight the whole paragraph (incl. paragraph mark) cyan so the
        # <w:pPr><w:rPr> also ends up with <w:highlight w:val="cyan"/>.
        $p.Range.Font.HighlightColorIndex = 3

        # Split point: right after "...stuck at si" / before "de of platforms..."
        $splitOffset = $full.IndexOf("de of platforms and ground")
        $bmPos = $p.Range.Start + $splitOffset
        $d.Bookmarks.Add("_GoBack", $d.Range($bmPos, $bmPos))

        break
    }
}
